# Replaced echo mechanism with simplex keep_alive
#
# - Row 57 ("Only send EchoReq as keep-alive when no incoming comm.") is
#   marked Rejected, with a note explaining it was replaced by the simplex
#   keep-alive mechanism.
# - Two new todo rows are appended (71, 72) describing the follow-up work.
# - The sheet's scroll/selection view state is moved further down the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 57: Status (C57) Ongoing -> Rejected
$ws.Range("C57").Value = "Rejected"

# Append new todo rows 71 & 72 (Status column = Open) before setting D57 so
# new shared-string entries are created in the same order as upstream.
$ws.Range("B71").Value = "Add SW power down"
$ws.Range("C71").Value = "Open"

$ws.Range("B72").Value = "Redo server side design diagram in draw.io"
$ws.Range("C72").Value = "Open"

# Row 57: add explanatory Notes (D57), inherits column D's wrap-text style
$ws.Range("D57").Value = "Replaced by simplex keep alive mechanism"

# Update the sheet view: scrolled down and new active selection
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D60").Select()
